# Edit script: "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Replaces the worker arrears (EC) data table with an updated dataset and
# shifts the signature footer rows down to make room for the extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update "Valor Mora" total (E11) ---------------------------------
$ws.Range("E11").Value = 1500987

# --- 2. New data table (rows 16..46) -------------------------------------
$data = @(
    @("CC","73204062","MARLON RAFAEL PARRA ORTIZ","2402","52000","1300000"),
    @("CC","73204062","MARLON RAFAEL PARRA ORTIZ","2401","26000","1300000"),
    @("CC","1049828499","EDINSON JAVIER CASTAÑO VICTOR","2403","46400","1160000"),
    @("CC","1049828499","EDINSON JAVIER CASTAÑO VICTOR","2402","46400","1160000"),
    @("CC","1049828499","EDINSON JAVIER CASTAÑO VICTOR","2401","46400","1160000"),
    @("CC","1049828499","EDINSON JAVIER CASTAÑO VICTOR","2305","126400","1160000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2412","17013","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2411","46400","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2410","46400","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2406","52000","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2405","52000","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2311","46400","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2310","46400","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2309","46400","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2308","46400","1000000"),
    @("CC","1101874652","LUIS ANGEL BERRIO LUNA","2307","6187","1000000"),
    @("CC","1007901653","MACHANISM BATISTA ABELLO","2407","52000","908526"),
    @("CC","1007901653","MACHANISM BATISTA ABELLO","2406","52000","908526"),
    @("CC","1052731257","JORGE ARMANDO GUERRA VERGAÑO","2404","52000","1300000"),
    @("CC","1052731257","JORGE ARMANDO GUERRA VERGAÑO","2403","52000","1300000"),
    @("CC","1052731257","JORGE ARMANDO GUERRA VERGAÑO","2402","52000","1300000"),
    @("CC","1047518591","DERWIN JOSE PEREZ LOPEZ","2407","52000","1300000"),
    @("CC","1047518591","DERWIN JOSE PEREZ LOPEZ","2406","52000","1300000"),
    @("CC","1047518591","DERWIN JOSE PEREZ LOPEZ","2405","52000","1300000"),
    @("CC","1193524781","OMAR YESITH JERONIMO BRAVO","2407","52000","1300000"),
    @("CC","1193524781","OMAR YESITH JERONIMO BRAVO","2406","52000","1300000"),
    @("CC","1193524781","OMAR YESITH JERONIMO BRAVO","2405","52000","1300000"),
    @("CC","1063174273","LUIS ALFREDO LOPEZ NUÑEZ","2404","52000","1300000"),
    @("CC","1063174273","LUIS ALFREDO LOPEZ NUÑEZ","2403","52000","1300000"),
    @("CC","1065122249","GERARDO JUNIOR CONTRERAS DIAZ","2402","29387","1160000"),
    @("CC","1065122249","GERARDO JUNIOR CONTRERAS DIAZ","2401","46400","1160000")
)

# --- 3. Move the footer (signature) block down from rows 48-49 to 51-52 --
# Copy formatting first (values/merges are reapplied explicitly afterwards)
$ws.Range("B48:C49").Copy()
$ws.Range("B51").PasteSpecial(-4122)
$ws.Range("H48:J49").Copy()
$ws.Range("H51").PasteSpecial(-4122)

$ws.Range("B51").Value = "___________________________________"
$ws.Range("H51").Value = "___________________________________"
$ws.Range("B52").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H52").Value = "FIRMA DEL REPRESENTANTE LEGAL"

$ws.Range("B51:C51").Merge()
$ws.Range("H51:J51").Merge()
$ws.Range("B52:C52").Merge()
$ws.Range("H52:J52").Merge()

# Remove the old footer cells/merges now that they have moved down
$ws.Range("B48:C48").UnMerge()
$ws.Range("H48:J48").UnMerge()
$ws.Range("B49:C49").UnMerge()
$ws.Range("H49:J49").UnMerge()
$ws.Range("B48:J49").Clear()

# --- 4. Write out the new data rows (16..46) ------------------------------
# Row 16 carries the "interior" row style; row 43 (old last row) carries the
# special bottom-border style used for the final row of the table. Move that
# special style onto the new true last row (46) BEFORE the row 43 source is
# overwritten by the interior-style loop below.
$ws.Range("B43:J43").Copy()
$ws.Range("B46").PasteSpecial(-4122)

$r = 16
foreach ($row in $data) {
    if ($r -ne 46) {
        $ws.Range("B16:J16").Copy()
        $ws.Cells.Item($r, 2).PasteSpecial(-4122)
    }

    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = [double]$row[4]
    $ws.Cells.Item($r, 7).Value = [double]$row[5]
    $r = $r + 1
}

Write-Output "done"
